$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 42-43: the access-time / runtime values now duplicate already-used
# text labels ("3/25/24 18:53" / "00:00:00"), so the cells become text.
$ws.Range("A42").Value = "3/25/24 18:53"
$ws.Range("C42").Value = "00:00:00"
$ws.Range("A43").Value = "3/25/24 18:53"
$ws.Range("C43").Value = "00:00:00"

# New rows 44-47 - copy the number formatting (styles) down from the last
# existing data row first so the new cells share the same style indices,
# then overwrite with the new row's values.
$ws.Cells.Item(41, 1).Copy($ws.Cells.Item(44, 1))
$ws.Cells.Item(41, 3).Copy($ws.Cells.Item(44, 3))
$ws.Cells.Item(41, 1).Copy($ws.Cells.Item(45, 1))
$ws.Cells.Item(41, 3).Copy($ws.Cells.Item(45, 3))
$ws.Cells.Item(41, 1).Copy($ws.Cells.Item(46, 1))
$ws.Cells.Item(41, 3).Copy($ws.Cells.Item(46, 3))
$ws.Cells.Item(41, 1).Copy($ws.Cells.Item(47, 1))
$ws.Cells.Item(41, 3).Copy($ws.Cells.Item(47, 3))

# Row 44: test1 sample, sharing the "4/3/24 20:31" / "00:00:00" text labels.
$ws.Range("A44").Value = "4/3/24 20:31"
$ws.Range("B44").Value = "test1"
$ws.Range("C44").Value = "00:00:00"

# Row 45: test2 sample, same text labels.
$ws.Range("A45").Value = "4/3/24 20:31"
$ws.Range("B45").Value = "test2"
$ws.Range("C45").Value = "00:00:00"

# Row 46: test1 sample, real numeric date serial + elapsed-time values.
$ws.Range("A46").Value = 45385.855223324834
$ws.Range("B46").Value = "test1"
$ws.Range("C46").Value = 0.000000000000474537

# Row 47: test2 sample, real numeric date serial + elapsed-time values.
$ws.Range("A47").Value = 45385.855223324834
$ws.Range("B47").Value = "test2"
$ws.Range("C47").Value = 0.000000000000474537
